$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing "Regex" section (old row 47),
# pushing the Regex section and everything after it down by two rows.
$ws.Rows("47:48").Insert()

# The insert copies formatting down from the row above; start clean.
$ws.Rows("47:48").ClearFormats()
$ws.Rows("47:48").ClearContents()

# New "CRMInfo" section header (matches style of other section headers,
# e.g. the "MailBoxes" header in A41).
$ws.Range("A41").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A47").Value = "CRMInfo"

# New "DocumentType" / "LTRCUST" parameter row (matches style of a normal
# name/value row, e.g. row 42).
$ws.Range("A42:B42").Copy()
$ws.Range("A48:B48").PasteSpecial(-4122)
$ws.Range("A48").Value = "DocumentType"
$ws.Range("B48").Value = "LTRCUST"

# Resize the worksheet table to include the two new rows.
$tbl = $ws.ListObjects.Item(1)
$null = $tbl.Resize($ws.Range("A1:C57"))

# Update the view / selection to match the saved workbook state.
$null = $ws.Range("B49").Select()
